$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-23 Wednesday", "2025-04-24 Thursday"),
    @("83÷8=", "84÷7="),
    @("19÷7=", "95÷9="),
    @("92÷3=", "25÷9="),
    @("83÷6=", "63÷9="),
    @("17÷2=", "98÷3="),
    @("72÷6=", "35÷2="),
    @("93÷3=", "29÷9="),
    @("80÷4=", "31÷2="),
    @("81÷7=", "17÷6="),
    @("76÷5=", "36÷4="),
    @("47÷4=", "34÷8="),
    @("11÷4=", "82÷9="),
    @("88÷8=", "52÷9="),
    @("84÷5=", "36÷6="),
    @("16÷4=", "66÷7="),
    @("70÷9=", "71÷7="),
    @("48÷6=", "87÷5="),
    @("57÷9=", "68÷2="),
    @("17÷3=", "68÷5="),
    @("54÷2=", "74÷4="),
    @("96÷8=", "53÷8="),
    @("38÷4=", "56÷2="),
    @("55÷9=", "80÷8="),
    @("54÷9=", "75÷8="),
    @("89÷3=", "71÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
